$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be auto-detected as numbers,
# so they stay text cells like the rest of the Price column.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '67.778.95'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.424.73'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '553.77'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '161.02'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('E9').Value = '  +7.51%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = '4.79'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  -2.29%  '
$ws.Range('D13').Value = '67.741.77'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').Value = '23.09'
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('D16').Value = '10.33'
$ws.Range('E16').Value = '  -3.90%  '
$ws.Range('D17').Value = '335.86'
$ws.Range('D18').Value = '6.84'
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '66.42'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('D24').Value = '8.09'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('D25').Value = '0.0₃0813'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '423.30'
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('D31').Value = '160.74'
$ws.Range('E31').Value = '  +2.61%  '
$ws.Range('D32').Value = '18.93'
$ws.Range('E32').Value = '  -0.46%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -5.85%  '
$ws.Range('D36').Value = '0.294'
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('E37').Value = '  -3.10%  '
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('D40').Value = '2.03'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D42').Value = '128.99'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '0.555'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').Value = '0.0913'
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('E48').Value = '  -5.41%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0207'
$ws.Range('E49').Value = '  +4.90%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '16.62'
$ws.Range('E50').Value = '  -0.84%  '
$ws.Range('D51').Value = '4.77'
$ws.Range('E51').Value = '  -6.31%  '
